$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - candidate/member info
$ws.Range("A2").Value = "PmJ Tangará"
$ws.Range("B2").Value = "Baltazar Patricio Marinho De Figueiredo"
$ws.Range("C2").Value = 9
$ws.Range("E2").Value = "9ª PmJ Parnamirim"

# Row 6 - correções de fls.
$ws.Range("C6").Value = "fls. 65-69"

# Row 7
$ws.Range("C7").Value = "fls. 10"

# Row 9
$ws.Range("C9").Value = "fls. 10"

# Row 11
$ws.Range("C11").Value = "fls. 3-4"

# Row 12
$ws.Range("C12").Value = "fls. 3-4"

# Row 13
$ws.Range("C13").Value = "fls. 5"

# Row 14
$ws.Range("C14").Value = "fls. 5"

# Row 15
$ws.Range("C15").Value = "fls. 5"

# Row 16
$ws.Range("C16").Value = "fls. 5"

# Row 17
$ws.Range("C17").Value = "fls. 6"

# Row 18
$ws.Range("C18").Value = "fls. 6"

# Row 19
$ws.Range("C19").Value = "fls. 6"

# Row 20
$ws.Range("C20").Value = "fls. 8"

# Row 21
$ws.Range("C21").Value = "fls. 6"

# Row 22
$ws.Range("C22").Value = "fls. 6"

# Row 23
$ws.Range("C23").Value = "fls. 6"

# Row 24
$ws.Range("C24").Value = "fls. 6"
